$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(3).TextFrame.TextRange
# First trim the text range down to the content of the existing first run so
# that the writer keeps (and simply extends) that run instead of creating a
# brand-new one -- this avoids an unwanted lang="en-US" stamp and mirrors the
# "consolidate text runs" behaviour from the diff.
$tr.Text = "Followed"
$tr.Text = "Followed by a picture"
